$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.018094897270203
$ws.Range("B1").Value = 1.482462763786316
$ws.Range("C1").Value = 3.038525581359863
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.458145141601562
